$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6482.381
$ws.Range("J17").Value = 6731.5
$ws.Range("L17").Value = 20194.5
$ws.Range("N17").Value = -20530.5
$ws.Range("H45").Value = 1995
$ws.Range("I45").Value = 742.4545000000001
$ws.Range("K45").Value = 2227.3635
$ws.Range("M45").Value = -2035.3635
$ws.Range("H69").Value = 9442.799999999999
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9442.799999999999
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 28328.4
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -30076.4
$ws.Range("H72").Value = 9442.799999999999
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9442.799999999999
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 84985.2
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -93721.2
$ws.Range("H86").Value = 6028.7144
$ws.Range("I86").Value = 5100.5
$ws.Range("J86").Value = 6400
$ws.Range("K86").Value = 5100.5
$ws.Range("L86").Value = 6400
$ws.Range("M86").Value = -3977.5
$ws.Range("N86").Value = -8646
$ws.Range("H89").Value = 6028.7144
$ws.Range("I89").Value = 5100.5
$ws.Range("J89").Value = 6400
$ws.Range("K89").Value = 25502.5
$ws.Range("L89").Value = 32000
$ws.Range("M89").Value = -19886.5
$ws.Range("N89").Value = -43232
$ws.Range("H112").Value = 1944
$ws.Range("I112").Value = 1240
$ws.Range("J112").Value = 2026.8235
$ws.Range("K112").Value = 3720
$ws.Range("L112").Value = 6080.470499999999
$ws.Range("M112").Value = -2612
$ws.Range("N112").Value = -8296.470499999999
$ws.Range("H121").Value = 812.75
$ws.Range("J121").Value = 812.75
$ws.Range("L121").Value = 2438.25
$ws.Range("N121").Value = -5932.25
$ws.Range("H138").Value = 2841.7615
$ws.Range("I138").Value = 1577.826
$ws.Range("J138").Value = 3289
$ws.Range("K138").Value = 4733.478
$ws.Range("L138").Value = 9867
$ws.Range("M138").Value = 406.5219999999999
$ws.Range("N138").Value = -20147

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 100001780
$ws.Range("I45").Value = 111112860
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 111112860
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -111112483
$ws.Range("N45").Value = -2754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 60749.25
$ws.Range("J132").Value = 60749.25
$ws.Range("L132").Value = 60749.25
$ws.Range("N132").Value = -70869.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1265.375
$ws.Range("I16").Value = 655.0625
$ws.Range("K16").Value = 655.0625
$ws.Range("M16").Value = -368.0625
$ws.Range("H113").Value = 1265.375
$ws.Range("I113").Value = 655.0625
$ws.Range("K113").Value = 655.0625
$ws.Range("M113").Value = 1514.9375
$ws.Range("H132").Value = 3498.8333
$ws.Range("I132").Value = 2972
$ws.Range("K132").Value = 8916
$ws.Range("M132").Value = -6386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2540.2
$ws.Range("I5").Value = 674
$ws.Range("J5").Value = 10005
$ws.Range("K5").Value = 2022
$ws.Range("L5").Value = 30015
$ws.Range("M5").Value = -1910
$ws.Range("N5").Value = -30239
$ws.Range("H102").Value = 7000
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H122").Value = 7660733
$ws.Range("I122").Value = 5128869
$ws.Range("J122").Value = 8405399
$ws.Range("K122").Value = 46159821
$ws.Range("L122").Value = 75648591
$ws.Range("M122").Value = -46157371
$ws.Range("N122").Value = -75653491
$ws.Range("H128").Value = 125988.5
$ws.Range("I128").Value = 125988.5
$ws.Range("K128").Value = 377965.5
$ws.Range("M128").Value = -372985.5
$ws.Range("H135").Value = 2540.2
$ws.Range("I135").Value = 674
$ws.Range("J135").Value = 10005
$ws.Range("K135").Value = 6066
$ws.Range("L135").Value = 90045
$ws.Range("M135").Value = -3531
$ws.Range("N135").Value = -95115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 359809.78
$ws.Range("I80").Value = 387179.78
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 387179.78
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -386181.78
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 359809.78
$ws.Range("I83").Value = 387179.78
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 1935898.9
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -1930906.9
$ws.Range("N83").Value = -29984
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""
$ws.Range("H102").Value = 1959.4
$ws.Range("I102").Value = 1149.8636
$ws.Range("K102").Value = 1149.8636
$ws.Range("M102").Value = 472.1364000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6442
$ws.Range("I7").Value = 4415.3335
$ws.Range("J7").Value = 9144.223
$ws.Range("K7").Value = 4415.3335
$ws.Range("L7").Value = 9144.223
$ws.Range("M7").Value = -4303.3335
$ws.Range("N7").Value = -9368.223
$ws.Range("H16").Value = 2583.8125
$ws.Range("I16").Value = 605.7
$ws.Range("K16").Value = 605.7
$ws.Range("M16").Value = -435.7
$ws.Range("H126").Value = 6442
$ws.Range("I126").Value = 4415.3335
$ws.Range("J126").Value = 9144.223
$ws.Range("K126").Value = 13246.0005
$ws.Range("L126").Value = 27432.669
$ws.Range("M126").Value = -10776.0005
$ws.Range("N126").Value = -32372.669
$ws.Range("H132").Value = 5183.321
$ws.Range("I132").Value = 5382.975
$ws.Range("J132").Value = 4569
$ws.Range("K132").Value = 16148.925
$ws.Range("L132").Value = 13707
$ws.Range("M132").Value = -13618.925
$ws.Range("N132").Value = -18767

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1350.5
$ws.Range("J4").Value = 1350.5
$ws.Range("L4").Value = 1350.5
$ws.Range("N4").Value = -1576.5
$ws.Range("H5").Value = 23287714
$ws.Range("J5").Value = 27167334
$ws.Range("L5").Value = 27167334
$ws.Range("N5").Value = -27167558
$ws.Range("H46").Value = 69619.336
$ws.Range("J46").Value = 69619.336
$ws.Range("L46").Value = 69619.336
$ws.Range("N46").Value = -70081.336
$ws.Range("H134").Value = 69619.336
$ws.Range("J134").Value = 69619.336
$ws.Range("L134").Value = 208858.008
$ws.Range("N134").Value = -213928.008
$ws.Range("H136").Value = 3125.5881
$ws.Range("I136").Value = 1975.1578
$ws.Range("K136").Value = 5925.4734
$ws.Range("M136").Value = -3375.4734
